$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.221.60'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.70%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.656.77'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('E5').Value = '  -0.80%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.5234'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.99%  '
$ws.Range('E7').Value = '  -0.60%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2667'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.19%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06365'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  -1.42%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07738'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.46%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '4.592'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.650.04'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.36%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.884.35'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.87%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.5643'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.35%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0₅8276'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.04%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '65.44'
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '26.226.80'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.77%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.692'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.42%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '192.32'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -3.52%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '10.39'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.68%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.020'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('E24').Value = '  -0.61%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '143.33'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -2.29%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.1200'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -2.63%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.275'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.32%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.96'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.26%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.499'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.08%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.05628'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -4.74%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.279'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.54%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.499'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.62%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.363'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.42%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.581'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.17%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.807'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -1.31%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.9459'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -2.55%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.411'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.03%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.5743'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.53%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01596'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.11%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '5.903'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.03%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.573'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.15%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.8451'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.31%  '
$ws.Range('E43').Value = '  -0.61%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.023.84'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -5.08%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '101.57'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.24%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.795.16'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.85%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '58.44'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.11%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0₈106'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.76%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.05314'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.85%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.4349'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.37%  '
